# Refresh the scraped "cryptos" price/volume snapshot (GitHub Actions job).
#
# The sheet stores Price (col D) and Volume(1h) (col E) as literal text
# (inline strings in the OOXML), not numbers - e.g. "1.000" must stay the
# 5-character string "1.000", not become the number 1. Plain
# `$ws.Range(...).Value = "1.000"` lets Excel's COM layer auto-detect a
# numeric-looking string and silently coerce it to a real number (dropping
# trailing zeros / reformatting), so any target value that parses as a bare
# decimal number is written as Text (NumberFormat "@") first. Values that
# are not valid bare numbers (two dots, like "29.179.88") and the Volume
# strings (they carry padding spaces / "%") already round-trip as text and
# need no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 47/48 swapped rank: BabyDogeCoin and Mantle traded places ---
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

# --- Cells whose new value parses as a plain decimal number and must be ---
# --- forced to Text so Excel doesn't rewrite it as a numeric literal.   ---
$riskyCells = @(
    @{ Cell = 'D4';  Value = '0.9983' }
    @{ Cell = 'D5';  Value = '241.36' }
    @{ Cell = 'D6';  Value = '0.6190' }
    @{ Cell = 'D7';  Value = '1.001' }
    @{ Cell = 'D8';  Value = '0.07349' }
    @{ Cell = 'D9';  Value = '0.2899' }
    @{ Cell = 'D10'; Value = '23.02' }
    @{ Cell = 'D11'; Value = '0.07669' }
    @{ Cell = 'D13'; Value = '4.971' }
    @{ Cell = 'D14'; Value = '0.6616' }
    @{ Cell = 'D15'; Value = '82.17' }
    @{ Cell = 'D16'; Value = '0.000008929' }
    @{ Cell = 'D17'; Value = '5.838' }
    @{ Cell = 'D20'; Value = '237.15' }
    @{ Cell = 'D21'; Value = '12.42' }
    @{ Cell = 'D22'; Value = '1.000' }
    @{ Cell = 'D23'; Value = '7.188' }
    @{ Cell = 'D24'; Value = '1.001' }
    @{ Cell = 'D25'; Value = '158.34' }
    @{ Cell = 'D26'; Value = '0.1413' }
    @{ Cell = 'D27'; Value = '8.436' }
    @{ Cell = 'D28'; Value = '17.64' }
    @{ Cell = 'D29'; Value = '1.486' }
    @{ Cell = 'D30'; Value = '0.05569' }
    @{ Cell = 'D31'; Value = '4.092' }
    @{ Cell = 'D32'; Value = '4.097' }
    @{ Cell = 'D33'; Value = '1.206' }
    @{ Cell = 'D34'; Value = '1.828' }
    @{ Cell = 'D35'; Value = '0.7334' }
    @{ Cell = 'D37'; Value = '2.620' }
    @{ Cell = 'D38'; Value = '2.843' }
    @{ Cell = 'D40'; Value = '0.01760' }
    @{ Cell = 'D41'; Value = '6.308' }
    @{ Cell = 'D42'; Value = '0.9213' }
    @{ Cell = 'D43'; Value = '1.000' }
    @{ Cell = 'D44'; Value = '101.58' }
    @{ Cell = 'D46'; Value = '64.76' }
    @{ Cell = 'D47'; Value = '0.5080' }
    @{ Cell = 'D48'; Value = '0.00000000117' }
    @{ Cell = 'D49'; Value = '0.4013' }
    @{ Cell = 'D50'; Value = '9.037' }
    @{ Cell = 'D51'; Value = '0.05773' }
)

foreach ($item in $riskyCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = '@'
    $rng.Value = $item.Value
}

# --- Cells that already round-trip safely as text (multi-dot prices, ---
# --- and every Volume(1h) percentage string).                        ---
$safeCells = @(
    @{ Cell = 'D2';  Value = '29.179.88' }
    @{ Cell = 'E2';  Value = '  +0.22%  ' }
    @{ Cell = 'D3';  Value = '1.825.17' }
    @{ Cell = 'E3';  Value = '  -0.39%  ' }
    @{ Cell = 'E4';  Value = '  -0.04%  ' }
    @{ Cell = 'E5';  Value = '  -0.69%  ' }
    @{ Cell = 'E6';  Value = '  -1.44%  ' }
    @{ Cell = 'E7';  Value = '  +0.00%  ' }
    @{ Cell = 'E8';  Value = '  -2.41%  ' }
    @{ Cell = 'E9';  Value = '  -1.09%  ' }
    @{ Cell = 'E10'; Value = '  -0.85%  ' }
    @{ Cell = 'E11'; Value = '  -0.20%  ' }
    @{ Cell = 'D12'; Value = '1.825.71' }
    @{ Cell = 'E12'; Value = '  -0.35%  ' }
    @{ Cell = 'E13'; Value = '  -1.16%  ' }
    @{ Cell = 'E14'; Value = '  -1.06%  ' }
    @{ Cell = 'E15'; Value = '  -1.06%  ' }
    @{ Cell = 'E16'; Value = '  -4.74%  ' }
    @{ Cell = 'E17'; Value = '  -2.66%  ' }
    @{ Cell = 'D18'; Value = '29.159.99' }
    @{ Cell = 'E18'; Value = '  +0.19%  ' }
    @{ Cell = 'D19'; Value = '2.066.46' }
    @{ Cell = 'E19'; Value = '  -0.33%  ' }
    @{ Cell = 'E20'; Value = '  +6.23%  ' }
    @{ Cell = 'E21'; Value = '  -1.44%  ' }
    @{ Cell = 'E22'; Value = '  -0.15%  ' }
    @{ Cell = 'E23'; Value = '  +0.62%  ' }
    @{ Cell = 'E25'; Value = '  -1.09%  ' }
    @{ Cell = 'E26'; Value = '  +1.07%  ' }
    @{ Cell = 'E27'; Value = '  -0.77%  ' }
    @{ Cell = 'E28'; Value = '  -1.48%  ' }
    @{ Cell = 'E29'; Value = '  -0.71%  ' }
    @{ Cell = 'E30'; Value = '  -4.08%  ' }
    @{ Cell = 'E31'; Value = '  -0.77%  ' }
    @{ Cell = 'E32'; Value = '  -1.67%  ' }
    @{ Cell = 'E33'; Value = '  +0.14%  ' }
    @{ Cell = 'E34'; Value = '  -0.42%  ' }
    @{ Cell = 'E35'; Value = '  -1.06%  ' }
    @{ Cell = 'E36'; Value = '  -0.71%  ' }
    @{ Cell = 'E37'; Value = '  -1.80%  ' }
    @{ Cell = 'E38'; Value = '  +2.78%  ' }
    @{ Cell = 'D39'; Value = '1.216.35' }
    @{ Cell = 'E39'; Value = '  -0.80%  ' }
    @{ Cell = 'E40'; Value = '  -1.10%  ' }
    @{ Cell = 'E41'; Value = '  -2.84%  ' }
    @{ Cell = 'E42'; Value = '  +3.08%  ' }
    @{ Cell = 'E43'; Value = '  +0.01%  ' }
    @{ Cell = 'E44'; Value = '  -0.36%  ' }
    @{ Cell = 'D45'; Value = '1.971.47' }
    @{ Cell = 'E45'; Value = '  -0.26%  ' }
    @{ Cell = 'E46'; Value = '  -1.75%  ' }
    @{ Cell = 'E47'; Value = '  -0.21%  ' }
    @{ Cell = 'E48'; Value = '  -2.65%  ' }
    @{ Cell = 'E49'; Value = '  -1.39%  ' }
    @{ Cell = 'E50'; Value = '  +0.32%  ' }
    @{ Cell = 'E51'; Value = '  -0.84%  ' }
)

foreach ($item in $safeCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
